$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026490568201804
$ws.Range("D2").Value = 1.027900252038954
$ws.Range("E2").Value = 1.035331548549453
$ws.Range("F2").Value = 1.043065498449459
$ws.Range("I2").Value = 1.029911313145199
$ws.Range("J2").Value = 1.031653717339011
$ws.Range("K2").Value = 1.030718634338915
$ws.Range("L2").Value = 1.038128458882054
$ws.Range("M2").Value = 1.045840415014454
$ws.Range("N2").Value = 1.01457077155899
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027358797952103
$ws.Range("D3").Value = 1.028641028618805
$ws.Range("E3").Value = 1.036137480706631
$ws.Range("F3").Value = 1.044033996961413
$ws.Range("I3").Value = 1.029981840884453
$ws.Range("J3").Value = 1.032162147637047
$ws.Range("K3").Value = 1.031267414437845
$ws.Range("L3").Value = 1.038743773010692
$ws.Range("M3").Value = 1.046619460540588
$ws.Range("N3").Value = 1.0147406697769
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027921102406145
$ws.Range("D4").Value = 1.029121116591508
$ws.Range("E4").Value = 1.036659867040349
$ws.Range("F4").Value = 1.044661921992077
$ws.Range("I4").Value = 1.030026042117505
$ws.Range("J4").Value = 1.032491004217705
$ws.Range("K4").Value = 1.03162260266247
$ws.Range("L4").Value = 1.039142164146712
$ws.Range("M4").Value = 1.047124177144633
$ws.Range("N4").Value = 1.014850515695606
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028157613873335
$ws.Range("D5").Value = 1.029323124872392
$ws.Range("E5").Value = 1.036879690696496
$ws.Range("F5").Value = 1.044926197273333
$ws.Range("I5").Value = 1.030044280173663
$ws.Range("J5").Value = 1.032629222752566
$ws.Range("K5").Value = 1.031771944227302
$ws.Range("L5").Value = 1.039309704317523
$ws.Range("M5").Value = 1.047336507428002
$ws.Range("N5").Value = 1.014896673021644
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028197332106577
$ws.Range("D6").Value = 1.029357053416581
$ws.Range("E6").Value = 1.036916612464911
$ws.Range("F6").Value = 1.044970587541318
$ws.Range("I6").Value = 1.03004732222494
$ws.Range("J6").Value = 1.032652428295729
$ws.Range("K6").Value = 1.031797020493591
$ws.Range("L6").Value = 1.039337838325889
$ws.Range("M6").Value = 1.04737216721831
$ws.Range("N6").Value = 1.014904421742566
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027924262219441
$ws.Range("D7").Value = 1.029123815133654
$ws.Range("E7").Value = 1.036662803502142
$ws.Range("F7").Value = 1.044665452091919
$ws.Range("I7").Value = 1.030026287168319
$ws.Range("J7").Value = 1.032492851230827
$ws.Range("K7").Value = 1.031624598093111
$ws.Range("L7").Value = 1.039144402606202
$ws.Range("M7").Value = 1.047127013734879
$ws.Range("N7").Value = 1.014851132538666
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026783885750828
$ws.Range("D8").Value = 1.02815044340443
$ws.Range("E8").Value = 1.035603730915951
$ws.Range("F8").Value = 1.043392548939361
$ws.Range("I8").Value = 1.029935445065723
$ws.Range("J8").Value = 1.031825570303083
$ws.Range("K8").Value = 1.030904077557257
$ws.Range("L8").Value = 1.038336356311933
$ws.Range("M8").Value = 1.046103567395207
$ws.Range("N8").Value = 1.01462820773649
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024778303412466
$ws.Range("D9").Value = 1.026441102638928
$ws.Range("E9").Value = 1.033744432826015
$ws.Range("F9").Value = 1.041159115208024
$ws.Range("I9").Value = 1.029764407070385
$ws.Range("J9").Value = 1.030648780084708
$ws.Range("K9").Value = 1.029635181355723
$ws.Range("L9").Value = 1.036914386152657
$ws.Range("M9").Value = 1.044304954628535
$ws.Range("N9").Value = 1.014234719492239
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023443960746389
$ws.Range("D10").Value = 1.025305585586528
$ws.Range("E10").Value = 1.032509654467342
$ws.Range("F10").Value = 1.039676702799082
$ws.Range("I10").Value = 1.029643047751488
$ws.Range("J10").Value = 1.029863682593621
$ws.Range("K10").Value = 1.028789834304936
$ws.Range("L10").Value = 1.035967770162659
$ws.Range("M10").Value = 1.043109216176805
$ws.Range("N10").Value = 1.01397197277894
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022866836801791
$ws.Range("D11").Value = 1.024814874370965
$ws.Range("E11").Value = 1.031976130138172
$ws.Range("F11").Value = 1.039036374737424
$ws.Range("I11").Value = 1.02958876598452
$ws.Range("J11").Value = 1.029523606475531
$ws.Range("K11").Value = 1.028423944429272
$ws.Range("L11").Value = 1.035558215230219
$ws.Range("M11").Value = 1.042592258324716
$ws.Range("N11").Value = 1.01385810625985
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022652567128421
$ws.Range("D12").Value = 1.02463275058978
$ws.Range("E12").Value = 1.031778128823551
$ws.Range("F12").Value = 1.038798765376878
$ws.Range("I12").Value = 1.029568343645714
$ws.Range("J12").Value = 1.029397269680555
$ws.Range("K12").Value = 1.028288060623182
$ws.Range("L12").Value = 1.035406140114066
$ws.Range("M12").Value = 1.042400359565393
$ws.Range("N12").Value = 1.013815797271415
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022698524159272
$ws.Range("D13").Value = 1.024671810024743
$ws.Range("E13").Value = 1.031820592905109
$ws.Range("F13").Value = 1.038849722640002
$ws.Range("I13").Value = 1.029572736047908
$ws.Range("J13").Value = 1.029424370119198
$ws.Range("K13").Value = 1.028317207054579
$ws.Range("L13").Value = 1.035438758376253
$ws.Range("M13").Value = 1.042441516937618
$ws.Range("N13").Value = 1.013824873314567
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022849123158522
$ws.Range("D14").Value = 1.024799816929033
$ws.Range("E14").Value = 1.031959759741219
$ws.Range("F14").Value = 1.0390167290293
$ws.Range("I14").Value = 1.029587083161434
$ws.Range("J14").Value = 1.029513163784282
$ws.Range("K14").Value = 1.028412711734655
$ws.Range("L14").Value = 1.035545643588818
$ws.Range("M14").Value = 1.042576393407746
$ws.Range("N14").Value = 1.013854609266497
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022941925446625
$ws.Range("D15").Value = 1.02487870589822
$ws.Range("E15").Value = 1.032045528046225
$ws.Range("F15").Value = 1.039119658631533
$ws.Range("I15").Value = 1.029595888500172
$ws.Range("J15").Value = 1.029567870224557
$ws.Range("K15").Value = 1.028471558537295
$ws.Range("L15").Value = 1.035611506002692
$ws.Range("M15").Value = 1.042659511516979
$ws.Range("N15").Value = 1.013872928740288
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023482276487587
$ws.Range("D16").Value = 1.025338173164997
$ws.Range("E16").Value = 1.032545086941549
$ws.Range("F16").Value = 1.039719232406533
$ws.Range("I16").Value = 1.02964661380527
$ws.Range("J16").Value = 1.029886249821731
$ws.Range("K16").Value = 1.028814120503042
$ws.Range("L16").Value = 1.035994958168593
$ws.Range("M16").Value = 1.043143542050357
$ws.Range("N16").Value = 1.013979527751872
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023821400842208
$ws.Range("D17").Value = 1.025626647125903
$ws.Range("E17").Value = 1.032858754194879
$ws.Range("F17").Value = 1.040095749717601
$ws.Range("I17").Value = 1.02967796906285
$ws.Range("J17").Value = 1.030085928658108
$ws.Range("K17").Value = 1.029029041914032
$ws.Range("L17").Value = 1.036235578491045
$ws.Range("M17").Value = 1.043447378135022
$ws.Range("N17").Value = 1.014046369255156
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024019269466381
$ws.Range("D18").Value = 1.025795003081308
$ws.Range("E18").Value = 1.03304182115281
$ws.Range("F18").Value = 1.040315516907217
$ws.Range("I18").Value = 1.029696090864418
$ws.Range("J18").Value = 1.030202385926709
$ws.Range("K18").Value = 1.029154416416529
$ws.Range("L18").Value = 1.036375960587379
$ws.Range("M18").Value = 1.043624678132809
$ws.Range("N18").Value = 1.014085347532771
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024086748218064
$ws.Range("D19").Value = 1.025852423999184
$ws.Range("E19").Value = 1.033104260887938
$ws.Range("F19").Value = 1.040390477402978
$ws.Range("I19").Value = 1.02970224155401
$ws.Range("J19").Value = 1.030242092770588
$ws.Range("K19").Value = 1.029197168298276
$ws.Range("L19").Value = 1.036423832698087
$ws.Range("M19").Value = 1.043685145958621
$ws.Range("N19").Value = 1.014098636537178
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023785009468791
$ws.Range("D20").Value = 1.025595686867433
$ws.Range("E20").Value = 1.032825089261066
$ws.Range("F20").Value = 1.040055337341448
$ws.Range("I20").Value = 1.029674622229357
$ws.Range("J20").Value = 1.030064506243113
$ws.Range("K20").Value = 1.029005981373521
$ws.Range("L20").Value = 1.03620975885363
$ws.Range("M20").Value = 1.043414771366448
$ws.Range("N20").Value = 1.014039198746974
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022804772738931
$ws.Range("D21").Value = 1.024762117981104
$ws.Range("E21").Value = 1.031918773777128
$ws.Range("F21").Value = 1.038967543254576
$ws.Range("I21").Value = 1.029582865457562
$ws.Range("J21").Value = 1.029487016733211
$ws.Range("K21").Value = 1.028384587314203
$ws.Range("L21").Value = 1.035514167116518
$ws.Range("M21").Value = 1.042536672258031
$ws.Range("N21").Value = 1.01384585314973
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022189037440667
$ws.Range("D22").Value = 1.024238878497201
$ws.Range("E22").Value = 1.031349941192854
$ws.Range("F22").Value = 1.038284975931734
$ws.Range("I22").Value = 1.029523672237197
$ws.Range("J22").Value = 1.02912382618273
$ws.Range("K22").Value = 1.027994031620092
$ws.Range("L22").Value = 1.035077121466067
$ws.Range("M22").Value = 1.041985285432282
$ws.Range("N22").Value = 1.013724208937868
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022515395069813
$ws.Range("D23").Value = 1.024516175715823
$ws.Range("E23").Value = 1.0316513943876
$ws.Range("F23").Value = 1.03864668717535
$ws.Range("I23").Value = 1.029555193865274
$ws.Range("J23").Value = 1.029316369439392
$ws.Range("K23").Value = 1.028201058933474
$ws.Range("L23").Value = 1.035308778706321
$ws.Range("M23").Value = 1.042277518323004
$ws.Range("N23").Value = 1.01378870228495
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023801452969087
$ws.Range("D24").Value = 1.025609676185406
$ws.Range("E24").Value = 1.0328403006554
$ws.Range("F24").Value = 1.040073597488284
$ws.Range("I24").Value = 1.029676135036005
$ws.Range("J24").Value = 1.030074186147601
$ws.Range("K24").Value = 1.029016401394558
$ws.Range("L24").Value = 1.036221425536631
$ws.Range("M24").Value = 1.043429504722173
$ws.Range("N24").Value = 1.014042438819573
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025296322749414
$ws.Range("D25").Value = 1.026882302860412
$ws.Range("E25").Value = 1.034224275323598
$ws.Range("F25").Value = 1.041735365334634
$ws.Range("I25").Value = 1.029809919840894
$ws.Range("J25").Value = 1.030953114425114
$ws.Range("K25").Value = 1.029963124217217
$ws.Range("L25").Value = 1.037281764910125
$ws.Range("M25").Value = 1.044769357542342
$ws.Range("N25").Value = 1.014336521833749
